$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 14216.5
$ws.Range("I13").Value = 12575
$ws.Range("J13").Value = 17499.5
$ws.Range("K13").Value = 12575
$ws.Range("L13").Value = 17499.5
$ws.Range("M13").Value = -12406
$ws.Range("N13").Value = -17837.5
$ws.Range("H19").Value = 487.25
$ws.Range("I19").Value = 428.23077
$ws.Range("J19").Value = 538.4
$ws.Range("K19").Value = 428.23077
$ws.Range("L19").Value = 538.4
$ws.Range("M19").Value = -253.23077
$ws.Range("N19").Value = -888.4
$ws.Range("H31").Value = 14684
$ws.Range("I31").Value = 14684
$ws.Range("K31").Value = 44052
$ws.Range("M31").Value = -43822
$ws.Range("H58").Value = 2051.0625
$ws.Range("I58").Value = 982.8
$ws.Range("J58").Value = 2536.6365
$ws.Range("K58").Value = 2948.4
$ws.Range("L58").Value = 7609.9095
$ws.Range("M58").Value = -2798.4
$ws.Range("N58").Value = -7909.9095
$ws.Range("H94").Value = 14000
$ws.Range("I94").Value = 14000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 14000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -13549
$ws.Range("H112").Value = 2358.7715
$ws.Range("J112").Value = 2398.7354
$ws.Range("L112").Value = 7196.206200000001
$ws.Range("N112").Value = -9412.206200000001
$ws.Range("H121").Value = 1179.0588
$ws.Range("I121").Value = 798.3333
$ws.Range("J121").Value = 1260.6428
$ws.Range("K121").Value = 2394.9999
$ws.Range("L121").Value = 3781.9284
$ws.Range("M121").Value = -647.9998999999998
$ws.Range("N121").Value = -7275.928400000001
$ws.Range("H129").Value = 1065.0864
$ws.Range("I129").Value = 394
$ws.Range("J129").Value = 1118.7733
$ws.Range("K129").Value = 1182
$ws.Range("L129").Value = 3356.3199
$ws.Range("M129").Value = 3818
$ws.Range("N129").Value = -13356.3199
$ws.Range("H137").Value = 2337.0952
$ws.Range("I137").Value = 1562.0952
$ws.Range("J137").Value = 3112.0952
$ws.Range("K137").Value = 4686.2856
$ws.Range("L137").Value = 9336.285600000001
$ws.Range("M137").Value = -2136.2856
$ws.Range("N137").Value = -14436.2856
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H70").Value = 38750
$ws.Range("J70").Value = 38750
$ws.Range("L70").Value = 38750
$ws.Range("N70").Value = -39290
$ws.Range("H73").Value = 38750
$ws.Range("J73").Value = 38750
$ws.Range("L73").Value = 38750
$ws.Range("N73").Value = -40622
$ws.Range("H132").Value = 4919.3096
$ws.Range("I132").Value = 5081.0645
$ws.Range("J132").Value = 4463.4546
$ws.Range("K132").Value = 15243.1935
$ws.Range("L132").Value = 13390.3638
$ws.Range("M132").Value = -12713.1935
$ws.Range("N132").Value = -18450.3638
$ws.Range("H141").Value = 67431.664
$ws.Range("J141").Value = 67431.664
$ws.Range("L141").Value = 67431.664
$ws.Range("N141").Value = -77791.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 69597.5
$ws.Range("J138").Value = 69597.5
$ws.Range("L138").Value = 69597.5
$ws.Range("N138").Value = -79877.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1691.2603
$ws.Range("I31").Value = 2165.2
$ws.Range("K31").Value = 2165.2
$ws.Range("M31").Value = -1870.2
$ws.Range("H34").Value = 1691.2603
$ws.Range("I34").Value = 2165.2
$ws.Range("K34").Value = 2165.2
$ws.Range("M34").Value = -1963.2
$ws.Range("H99").Value = 1857.6471
$ws.Range("I99").Value = 1900
$ws.Range("J99").Value = 1660
$ws.Range("K99").Value = 1900
$ws.Range("L99").Value = 1660
$ws.Range("M99").Value = -402
$ws.Range("N99").Value = -4656
$ws.Range("H126").Value = 1857.6471
$ws.Range("I126").Value = 1900
$ws.Range("J126").Value = 1660
$ws.Range("K126").Value = 5700
$ws.Range("L126").Value = 4980
$ws.Range("M126").Value = -3230
$ws.Range("N126").Value = -9920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 222832.52
$ws.Range("I113").Value = 250569.3
$ws.Range("J113").Value = 938.2
$ws.Range("K113").Value = 751707.8999999999
$ws.Range("L113").Value = 2814.6
$ws.Range("M113").Value = -749537.8999999999
$ws.Range("N113").Value = -7154.6
$ws.Range("H122").Value = 680.35297
$ws.Range("I122").Value = 613.5
$ws.Range("J122").Value = 840.8
$ws.Range("K122").Value = 5521.5
$ws.Range("L122").Value = 7567.2
$ws.Range("M122").Value = -3071.5
$ws.Range("N122").Value = -12467.2
$ws.Range("H131").Value = 3037.0386
$ws.Range("J131").Value = 3683.9268
$ws.Range("L131").Value = 11051.7804
$ws.Range("N131").Value = -21131.7804
$ws.Range("H134").Value = 3138.4075
$ws.Range("I134").Value = 2200.8948
$ws.Range("J134").Value = 5365
$ws.Range("K134").Value = 6602.6844
$ws.Range("L134").Value = 16095
$ws.Range("M134").Value = -1532.6844
$ws.Range("N134").Value = -26235
$ws.Range("H137").Value = 25645174
$ws.Range("I137").Value = 4287
$ws.Range("J137").Value = 33337440
$ws.Range("K137").Value = 12861
$ws.Range("L137").Value = 100012320
$ws.Range("M137").Value = -7761
$ws.Range("N137").Value = -100022520
$ws.Range("H139").Value = 1668
$ws.Range("I139").Value = 1538.15
$ws.Range("J139").Value = 2966.5
$ws.Range("K139").Value = 4614.450000000001
$ws.Range("L139").Value = 8899.5
$ws.Range("M139").Value = 525.5499999999993
$ws.Range("N139").Value = -19179.5
$ws.Range("H141").Value = 9338.333000000001
$ws.Range("I141").Value = 9806
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 29418
$ws.Range("L141").Value = 21000
$ws.Range("M141").Value = -24238
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13925778
$ws.Range("I11").Value = 17900142
$ws.Range("J11").Value = 15500
$ws.Range("K11").Value = 17900142
$ws.Range("L11").Value = 15500
$ws.Range("M11").Value = -17900003
$ws.Range("N11").Value = -15778
$ws.Range("H18").Value = 21333.334
$ws.Range("J18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("N18").Value = -7586
$ws.Range("H70").Value = 284834.94
$ws.Range("I70").Value = 377672.75
$ws.Range("K70").Value = 377672.75
$ws.Range("M70").Value = -377402.75
$ws.Range("H73").Value = 284834.94
$ws.Range("I73").Value = 377672.75
$ws.Range("K73").Value = 377672.75
$ws.Range("M73").Value = -376736.75
$ws.Range("H132").Value = 4201
$ws.Range("I132").Value = 3150
$ws.Range("K132").Value = 9450
$ws.Range("M132").Value = -6920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 3000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 3000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 3000
$ws.Range("N23").Value = -3460
$ws.Range("H108").Value = 26500
$ws.Range("J108").Value = 26500
$ws.Range("L108").Value = 26500
$ws.Range("N108").Value = -34180
$ws.Range("H132").Value = 3775.3044
$ws.Range("I132").Value = 3658.7144
$ws.Range("K132").Value = 10976.1432
$ws.Range("M132").Value = -8446.143199999999
$ws.Range("H133").Value = 35163
$ws.Range("J133").Value = 35163
$ws.Range("L133").Value = 35163
$ws.Range("N133").Value = -40223
$ws.Range("H141").Value = 39550
$ws.Range("I141").Value = 24325
$ws.Range("J141").Value = 70000
$ws.Range("K141").Value = 24325
$ws.Range("L141").Value = 70000
$ws.Range("M141").Value = -19145
$ws.Range("N141").Value = -80360
$ws.Range("M23").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 20050
$ws.Range("H40").Value = 16800
$ws.Range("I40").Value = 1333.3334
$ws.Range("J40").Value = 40000
$ws.Range("K40").Value = 1333.3334
$ws.Range("L40").Value = 40000
$ws.Range("M40").Value = -1184.3334
$ws.Range("N40").Value = -40298
$ws.Range("H62").Value = 4914.615
$ws.Range("I62").Value = 4825
$ws.Range("J62").Value = 4991.4287
$ws.Range("K62").Value = 4825
$ws.Range("L62").Value = 4991.4287
$ws.Range("M62").Value = -4201
$ws.Range("N62").Value = -6239.4287
$ws.Range("H65").Value = 4914.615
$ws.Range("I65").Value = 4825
$ws.Range("J65").Value = 4991.4287
$ws.Range("K65").Value = 24125
$ws.Range("L65").Value = 24957.1435
$ws.Range("M65").Value = -21005
$ws.Range("N65").Value = -31197.1435
$ws.Range("H123").Value = 23788.879
$ws.Range("J123").Value = 23788.879
$ws.Range("L123").Value = 23788.879
$ws.Range("N123").Value = -33588.879
$ws.Range("H126").Value = 8561.947
$ws.Range("I126").Value = 10338.733
$ws.Range("J126").Value = 1899
$ws.Range("K126").Value = 31016.199
$ws.Range("L126").Value = 5697
$ws.Range("M126").Value = -28546.199
$ws.Range("N126").Value = -10637
$ws.Range("H135").Value = 52452.332
$ws.Range("J135").Value = 52452.332
$ws.Range("L135").Value = 52452.332
$ws.Range("N135").Value = -62592.332
$ws.Range("H137").Value = 49750
$ws.Range("J137").Value = 49750
$ws.Range("L137").Value = 49750
$ws.Range("N137").Value = -59950
$ws.Range("H139").Value = 49866.668
$ws.Range("J139").Value = 49800
$ws.Range("L139").Value = 49800
$ws.Range("N139").Value = -60080
